$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# 1) Rename headers on existing sheets
$ws1.Range("B1").Value = "Weekly_PO_Qty"
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# 2) Add the new "PO Forecast" sheet at the end of the workbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "PO Forecast"

# 3) Header row values
$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

# Copy header style (bold + centered + border) from the Weekly Quantity sheet
$ws1.Range("A1:B1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

# 4) Data rows
$dates = @(45564.99999999999, 45634.99999999999, 45641.99999999999, 45648.99999999999, 45655.99999999999, 45662.99999999999, 45669.99999999999, 45676.99999999999, 45683.99999999999, 45690.99999999999)
$forecast = @(1616, 144, 0, 0, 0, 0, 0, 0, 0, 0)
$lower = @(1615.999049830921, 143.9991947640814, -3.200790739703065, -150.4007762868927, -297.6007621161742, -444.8007485919724, -592.000735529158, -739.2007221948073, -886.4007093378447, -1033.600696975725)
$upper = @(1615.999049972987, 143.9991949102708, -3.200790561011134, -150.4007760374359, -297.6007607875553, -444.8007451610983, -592.0007291775713, -739.2007132221477, -886.4006971375973, -1033.600680751216)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws3.Cells.Item($row, 1).Value = $dates[$i]
    $ws3.Cells.Item($row, 2).Value = $forecast[$i]
    $ws3.Cells.Item($row, 3).Value = $lower[$i]
    $ws3.Cells.Item($row, 4).Value = $upper[$i]
}

# Copy the date-format style from the Weekly Quantity sheet column A onto the new ds column
$ws1.Range("A2").Copy()
$ws3.Range("A2:A11").PasteSpecial(-4122)

# Match page margins used on the other sheets
$ws3.PageSetup.LeftMargin = 54
$ws3.PageSetup.RightMargin = 54
$ws3.PageSetup.TopMargin = 72
$ws3.PageSetup.BottomMargin = 72
$ws3.PageSetup.HeaderMargin = 36
$ws3.PageSetup.FooterMargin = 36

$wb.Worksheets.Item(1).Select()
